$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '40.846.19'
$ws.Range('E2').Value = '  +1.76%  '

$ws.Range('D3').Value = '2.241.20'
$ws.Range('E3').Value = '  +0.30%  '

$ws.Range('D4').Value = "'0.999"
$ws.Range('E4').Value = '  -0.17%  '

$ws.Range('D5').Value = "'303.30"
$ws.Range('E5').Value = '  +3.06%  '

$ws.Range('D6').Value = "'91.23"
$ws.Range('E6').Value = '  +3.84%  '

$ws.Range('D7').Value = "'0.519"
$ws.Range('E7').Value = '  +1.37%  '

$ws.Range('E8').Value = '  +0.01%  '

$ws.Range('D9').Value = "'0.480"
$ws.Range('E9').Value = '  +1.65%  '

$ws.Range('D10').Value = "'32.13"
$ws.Range('E10').Value = '  +6.02%  '

$ws.Range('D11').Value = "'52.79"
$ws.Range('E11').Value = '  +7.43%  '

$ws.Range('D12').Value = "'0.0790"
$ws.Range('E12').Value = '  +1.06%  '

$ws.Range('D13').Value = "'0.115"
$ws.Range('E13').Value = '  +3.01%  '

$ws.Range('D14').Value = "'6.54"
$ws.Range('E14').Value = '  +1.00%  '

$ws.Range('D15').Value = '2.581.74'
$ws.Range('E15').Value = '  +0.05%  '

$ws.Range('D16').Value = "'14.07"
$ws.Range('E16').Value = '  +1.70%  '

$ws.Range('D17').Value = '2.245.83'
$ws.Range('E17').Value = '  +0.57%  '

$ws.Range('D18').Value = "'0.749"
$ws.Range('E18').Value = '  +2.42%  '

$ws.Range('D19').Value = '40.729.53'
$ws.Range('E19').Value = '  +1.68%  '

$ws.Range('B20').Value = 'InternetComputer(DFINITY)'
$ws.Range('C20').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D20').Value = "'11.66"
$ws.Range('E20').Value = '  +2.68%  '

$ws.Range('B21').Value = 'ShibaInu'
$ws.Range('C21').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D21').Value = '0.0₃0900'
$ws.Range('E21').Value = '  +1.13%  '

$ws.Range('D22').Value = "'5.85"
$ws.Range('E22').Value = '  +0.18%  '

$ws.Range('D23').Value = "'66.18"
$ws.Range('E23').Value = '  +0.84%  '

$ws.Range('D24').Value = "'238.97"
$ws.Range('E24').Value = '  +0.94%  '

$ws.Range('D25').Value = "'2.55"
$ws.Range('E25').Value = '  +3.45%  '

$ws.Range('E26').Value = '  -0.10%  '

$ws.Range('D27').Value = "'1.86"
$ws.Range('E27').Value = '  +2.65%  '

$ws.Range('D28').Value = "'23.93"
$ws.Range('E28').Value = '  +4.91%  '

$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').Value = "'2.17"
$ws.Range('E29').Value = '  -0.86%  '

$ws.Range('B30').Value = 'Cosmos'
$ws.Range('C30').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D30').Value = "'9.46"
$ws.Range('E30').Value = '  +2.61%  '

$ws.Range('D31').Value = "'158.21"
$ws.Range('E31').Value = '  +1.76%  '

$ws.Range('D32').Value = "'33.14"
$ws.Range('E32').Value = '  +3.57%  '

$ws.Range('E33').Value = '  +0.00%  '

$ws.Range('D34').Value = "'5.08"
$ws.Range('E34').Value = '  +2.81%  '

$ws.Range('D35').Value = "'3.04"
$ws.Range('E35').Value = '  +5.16%  '

$ws.Range('D36').Value = "'0.0728"
$ws.Range('E36').Value = '  +1.46%  '

$ws.Range('E37').Value = '  -0.54%  '

$ws.Range('E38').Value = '  +6.83%  '

$ws.Range('E39').Value = '  +2.24%  '

$ws.Range('B40').Value = 'Celestia'
$ws.Range('C40').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D40').Value = "'16.25"
$ws.Range('E40').Value = '  +3.22%  '

$ws.Range('B41').Value = 'ARBITRUM'
$ws.Range('C41').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D41').Value = "'1.78"
$ws.Range('E41').Value = '  +5.46%  '

$ws.Range('D42').Value = "'3.88"
$ws.Range('E42').Value = '  +0.85%  '

$ws.Range('D43').Value = '2.096.16'
$ws.Range('E43').Value = '  -1.41%  '

$ws.Range('D44').Value = "'19.97"
$ws.Range('E44').Value = '  +10.49%  '

$ws.Range('D45').Value = "'0.0276"
$ws.Range('E45').Value = '  +3.14%  '

$ws.Range('D46').Value = "'10.14"
$ws.Range('E46').Value = '  +2.84%  '

$ws.Range('D47').Value = "'2.93"
$ws.Range('E47').Value = '  +9.57%  '

$ws.Range('D48').Value = "'1.83"
$ws.Range('E48').Value = '  -14.03%  '

$ws.Range('D49').Value = "'1.52"
$ws.Range('E49').Value = '  +2.88%  '

$ws.Range('D50').Value = '2.455.55'
$ws.Range('E50').Value = '  +0.45%  '

$ws.Range('E51').Value = '  +3.21%  '
